$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.33333
$ws.Range("I9").Value = 118.25
$ws.Range("J9").Value = 165.2
$ws.Range("K9").Value = 118.25
$ws.Range("L9").Value = 165.2
$ws.Range("M9").Value = 50.75
$ws.Range("N9").Value = -503.2

$ws.Range("H12").Value = 138.16667
$ws.Range("I12").Value = 156.66667
$ws.Range("J12").Value = 119.666664
$ws.Range("K12").Value = 156.66667
$ws.Range("L12").Value = 119.666664
$ws.Range("M12").Value = 13.33332999999999
$ws.Range("N12").Value = -459.666664

$ws.Range("H29").Value = 303
$ws.Range("I29").Value = 303
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 909
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -628
$ws.Range("M29").ClearContents()

$ws.Range("H33").Value = 256.03705
$ws.Range("I33").Value = 205.125
$ws.Range("J33").Value = 663.3333
$ws.Range("K33").Value = 205.125
$ws.Range("L33").Value = 663.3333
$ws.Range("M33").Value = 23.875
$ws.Range("N33").Value = -1121.3333

$ws.Range("H51").Value = 191359.8
$ws.Range("I51").Value = 400000
$ws.Range("J51").Value = 139199.75
$ws.Range("K51").Value = 400000
$ws.Range("L51").Value = 139199.75
$ws.Range("M51").Value = -399516
$ws.Range("N51").Value = -140167.75

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

$ws.Range("H100").Value = 670.55554
$ws.Range("I100").Value = 670.55554
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 670.55554
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -129.55554

$ws.Range("H107").Value = 353.77777
$ws.Range("I107").Value = 335.5
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 335.5
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1584.5
$ws.Range("N107").Value = -4340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 26000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 34000
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 34000
$ws.Range("M38").Value = -9533
$ws.Range("N38").Value = -34934

$ws.Range("H97").Value = 1676.4
$ws.Range("I97").Value = 1381.2858
$ws.Range("J97").Value = 1934.625
$ws.Range("K97").Value = 1381.2858
$ws.Range("L97").Value = 1934.625
$ws.Range("M97").Value = -885.2858000000001
$ws.Range("N97").Value = -2926.625

$ws.Range("H102").Value = 1049.5
$ws.Range("I102").Value = 699.5
$ws.Range("J102").Value = 1399.5
$ws.Range("K102").Value = 699.5
$ws.Range("L102").Value = 1399.5
$ws.Range("M102").Value = 922.5
$ws.Range("N102").Value = -4643.5

$ws.Range("H110").Value = 409.5
$ws.Range("I110").Value = 409.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 409.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1635.5

$ws.Range("H122").Value = 1966.1666
$ws.Range("I122").Value = 1504
$ws.Range("J122").Value = 2428.3333
$ws.Range("K122").Value = 4512
$ws.Range("L122").Value = 7284.999899999999
$ws.Range("M122").Value = -2062
$ws.Range("N122").Value = -12184.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4416.5
$ws.Range("I99").Value = 4416.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4416.5
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -2918.5
$ws.Range("M99").ClearContents()

$ws.Range("H110").Value = 72000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 72000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 72000
$ws.Range("N110").Value = -80180

$ws.Range("H111").Value = 72000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 72000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 72000
$ws.Range("N111").Value = -80180

$ws.Range("H112").Value = 103999.664
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 103999.664
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 103999.664
$ws.Range("N112").Value = -106953.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 703.45
$ws.Range("I22").Value = 738.8125
$ws.Range("J22").Value = 562
$ws.Range("K22").Value = 738.8125
$ws.Range("L22").Value = 562
$ws.Range("M22").Value = -388.8125
$ws.Range("N22").Value = -1262

$ws.Range("H132").Value = 1110.2858
$ws.Range("I132").Value = 987
$ws.Range("J132").Value = 1850
$ws.Range("K132").Value = 2961
$ws.Range("L132").Value = 5550
$ws.Range("M132").Value = -431
$ws.Range("N132").Value = -10610

$ws.Range("H134").Value = 1531.25
$ws.Range("I134").Value = 1607.1428
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4821.428400000001
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -2286.428400000001
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 999
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 999
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2997
$ws.Range("N31").Value = -3573

$ws.Range("H87").Value = 1199.6666
$ws.Range("I87").Value = 1199.6666
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 3598.9998
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -2350.9998

$ws.Range("H90").Value = 1199.6666
$ws.Range("I90").Value = 1199.6666
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 10796.9994
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -4556.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 235.28572
$ws.Range("I97").Value = 251.58333
$ws.Range("J97").Value = 137.5
$ws.Range("K97").Value = 251.58333
$ws.Range("L97").Value = 137.5
$ws.Range("M97").Value = 244.41667
$ws.Range("N97").Value = -1129.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H61").Value = 3289.8
$ws.Range("I61").Value = 3289.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3289.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3087.8

$ws.Range("H68").Value = 1550
$ws.Range("I68").Value = 1687.5
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1687.5
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -938.5
$ws.Range("N68").Value = -2498

$ws.Range("H71").Value = 1550
$ws.Range("I71").Value = 1687.5
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 8437.5
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -4693.5
$ws.Range("N71").Value = -12488

$ws.Range("H100").Value = 3999
$ws.Range("I100").Value = 3999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3999
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3458

$ws.Range("H113").Value = 3289.8
$ws.Range("I113").Value = 3289.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3289.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1119.8

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2499.4
$ws.Range("I132").Value = 2499.5
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 7498.5
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -4968.5
$ws.Range("N132").Value = -12557
